$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename renumbered "DB-xx" samples to "ST-xx"
$ws.Range("A22").Value = "ST-19"
$ws.Range("A15").Value = "ST-12"
$ws.Range("A16").Value = "ST-13"
$ws.Range("A23").Value = "ST-20"
$ws.Range("A24").Value = "ST-21"
$ws.Range("A25").Value = "ST-22"

# Fill in newly measured "Su from Ncor" values
$ws.Range("M5").Value = 50
$ws.Range("M8").Value = 190

# Merge header cell A1:A2 and center it
$ws.Range("A1:A2").Merge()
$ws.Range("A1:A2").HorizontalAlignment = -4108

$ws.Range("M6").Select()
